$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E12").Value = 3
$ws.Range("E15").Value = 1
$ws.Range("E21").Value = 1
$ws.Range("E22").Value = 25
$ws.Range("E23").Value = 33
$ws.Range("E24").Value = 46
$ws.Range("E25").Value = 46
$ws.Range("E26").Value = 48
